# Add a new worksheet "Deanza" after the existing sheets and populate it with
# the Deanza transfer-rate data, then make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Deanza"

# --- Write column A (group) data rows first so the new shared strings are
#     interned in the same order as the source workbook: asian, white, urm,
#     then the header labels group/transfer2/transfer4/year.
$ws.Range("A2").Value = "asian"
$ws.Range("A3").Value = "asian"
$ws.Range("A4").Value = "asian"
$ws.Range("A5").Value = "asian"
$ws.Range("A6").Value = "asian"

$ws.Range("A7").Value = "white"
$ws.Range("A8").Value = "white"
$ws.Range("A9").Value = "white"
$ws.Range("A10").Value = "white"
$ws.Range("A11").Value = "white"

$ws.Range("A12").Value = "urm"
$ws.Range("A13").Value = "urm"
$ws.Range("A14").Value = "urm"
$ws.Range("A15").Value = "urm"
$ws.Range("A16").Value = "urm"

# --- Header row, written in A, C, D, B order to match the shared-string order.
$ws.Range("A1").Value = "group"
$ws.Range("C1").Value = "transfer2"
$ws.Range("D1").Value = "transfer4"
$ws.Range("B1").Value = "year"

# --- Remaining numeric data (year / transfer2 / transfer4).
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.012
$ws.Range("D2").Value = 0.013

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0.03
$ws.Range("D3").Value = 0.035

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.009
$ws.Range("D4").Value = 0.057

$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.01
$ws.Range("D5").Value = 0.055

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 0.002
$ws.Range("D6").Value = 0.054

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = -0.018
$ws.Range("D7").Value = -0.016

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0.027
$ws.Range("D8").Value = -0.031

$ws.Range("B9").Value = 3
$ws.Range("C9").Value = -0.015
$ws.Range("D9").Value = -0.041

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = -0.02
$ws.Range("D10").Value = -0.014

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -0.04
$ws.Range("D11").Value = 0.012

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.028
$ws.Range("D12").Value = -0.008

$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 0.055
$ws.Range("D13").Value = -0.011

$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 0.064
$ws.Range("D14").Value = -0.03

$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 0.071
$ws.Range("D15").Value = -0.036

$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 0.096
$ws.Range("D16").Value = -0.046

# Match the selection/active cell seen in the source sheetView.
[void]$ws.Range("E22").Select()

[void]$wb.Worksheets.Item("Deanza").Activate()
